# "recursion - baby steps"
#
# 1. Sheet "08-BC": prefix the three requirement-name cells (B3:B5) with
#    their 1-based requirement number.
# 2. Sheet "provenance": refresh the run timestamp in B12.
# 3. Sheet "08-BC": widen column D (no data added, just a column-width tweak).

$wb = $excel.ActiveWorkbook

$bc = $wb.Worksheets.Item("08-BC")
$bc.Range("B3").Value = "1. assigned_region"
$bc.Range("B4").Value = "2. liquid_phase"
$bc.Range("B5").Value = "3. solid_phase"
$bc.Columns("D").ColumnWidth = 17.83333333333333

$prov = $wb.Worksheets.Item("provenance")
$prov.Range("B12").Value = 43435.51472035676
